# Updated cryptos list on Sun Oct  6 20:24:48 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper note: D-column "Price" cells are stored as text (original workbook
# uses inline strings for every value, including numeric-looking prices).
# Prefixing the value with a leading apostrophe forces Excel to keep it as
# text instead of re-interpreting it as a number, exactly matching the
# original cell type.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'62.697.74"
$ws.Range("E2").Value = "  +1.40%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'2.438.97"
$ws.Range("E3").Value = "  +1.60%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.12%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'566.51"
$ws.Range("E5").Value = "  +1.15%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'145.38"
$ws.Range("E6").Value = "  +2.38%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.10%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.533"
$ws.Range("E8").Value = "  +0.22%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +2.57%  "

# Row 10 - TRON
$ws.Range("E10").Value = "  +0.43%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'5.30"
$ws.Range("E11").Value = "  +1.46%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.356"
$ws.Range("E12").Value = "  +2.45%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "'26.93"
$ws.Range("E13").Value = "  +6.14%  "

# Row 14 - ShibaInu
$ws.Range("D14").Value = "'0.0000181"
$ws.Range("E14").Value = "  +6.32%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "'2.890.36"
$ws.Range("E15").Value = "  +2.12%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "'62.459.05"
$ws.Range("E16").Value = "  +1.16%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "'2.413.20"
$ws.Range("E17").Value = "  +0.64%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "'11.23"
$ws.Range("E18").Value = "  +0.79%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "'6.96"
$ws.Range("E19").Value = "  +2.85%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'324.22"
$ws.Range("E20").Value = "  +1.36%  "

# Row 21 - Polkadot
$ws.Range("D21").Value = "'4.17"
$ws.Range("E21").Value = "  +1.49%  "

# Row 22 - Dai
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.24%  "

# Row 23 - SuiNetwork
$ws.Range("D23").Value = "'1.84"
$ws.Range("E23").Value = "  +6.99%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'67.35"
$ws.Range("E24").Value = "  +2.83%  "

# Row 25 - Aptos
$ws.Range("D25").Value = "'8.57"
$ws.Range("E25").Value = "  -1.24%  "

# Row 26 - Bittensor
$ws.Range("D26").Value = "'580.61"
$ws.Range("E26").Value = "  +3.48%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  +9.32%  "

# Row 28 - WrappedeETH
$ws.Range("E28").Value = "  +1.53%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -1.44%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = "'8.42"
$ws.Range("E30").Value = "  +3.61%  "

# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  +4.65%  "

# Row 32 - Kaspa
$ws.Range("E32").Value = "  +0.00%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  -0.07%  "

# Row 34 - ImmutableX
$ws.Range("D34").Value = "'1.50"
$ws.Range("E34").Value = "  +0.65%  "

# Rows 35/36 swap coins: NEARProtocol now ranks above FirstDigitalUSD
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'4.85"
$ws.Range("E35").Value = "  +2.31%  "

$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.10%  "

# Row 37 - PolygonEcosystemToken
$ws.Range("D37").Value = "'0.383"
$ws.Range("E37").Value = "  +1.24%  "

# Row 38 - EthereumClassic
$ws.Range("D38").Value = "'18.79"
$ws.Range("E38").Value = "  +1.88%  "

# Row 39 - RenderToken
$ws.Range("D39").Value = "'5.39"
$ws.Range("E39").Value = "  -0.13%  "

# Row 40 - Monero
$ws.Range("D40").Value = "'148.02"
$ws.Range("E40").Value = "  -2.75%  "

# Row 41 - Stacks
$ws.Range("D41").Value = "'1.82"
$ws.Range("E41").Value = "  +2.59%  "

# Row 42 - USDe
$ws.Range("E42").Value = "  +0.19%  "

# Row 43 - dogwifhat
$ws.Range("D43").Value = "'2.45"
$ws.Range("E43").Value = "  +9.89%  "

# Row 44 - Aave
$ws.Range("D44").Value = "'148.25"
$ws.Range("E44").Value = "  +0.74%  "

# Row 45 - Filecoin
$ws.Range("E45").Value = "  +2.47%  "

# Row 46 - Hedera
$ws.Range("D46").Value = "'0.0536"
$ws.Range("E46").Value = "  +1.63%  "

# Row 47 - InjectiveProtocol
$ws.Range("D47").Value = "'20.53"
$ws.Range("E47").Value = "  +4.13%  "

# Row 48 - Mantle
$ws.Range("D48").Value = "'0.602"
$ws.Range("E48").Value = "  +2.68%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +3.44%  "

# Row 50 - Stellar
$ws.Range("D50").Value = "'0.0921"
$ws.Range("E50").Value = "  +0.56%  "

# Row 51 - BitgetToken
$ws.Range("E51").Value = "  +4.94%  "
